$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("List_ID")

$ws.Range("C22").Value = "24/07/2005"
$ws.Range("D22").Value = "nbrnneputiqtwys@gmail.com"
$ws.Range("E22").Value = "goevbULWSS5"
$ws.Range("F22").Value = "pass"
